$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - existing record, values updated and strategy/event columns refreshed
$ws.Range("A2").Value = -0.7413855195045471
$ws.Range("B2").Value = "Distress"
$ws.Range("C2").Value = 2.387104034423828
$ws.Range("D2").Value = "Talk"
$ws.Range("E2").Value = "Attention Deployment"

# Row 3 - new record
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "Love"
$ws.Range("C3").Value = 1.7927955389022827
$ws.Range("D3").Value = "Hello"
$ws.Range("E3").Value = "None"

# Row 4 - new record
$ws.Range("A4").Value = 0.8303518295288086
$ws.Range("B4").Value = "Love"
$ws.Range("C4").Value = 2.6735565662384033
$ws.Range("D4").Value = "Conversation"
$ws.Range("E4").Value = "None"

# Row 5 - new record
$ws.Range("A5").Value = 2.3478033542633057
$ws.Range("B5").Value = "Love"
$ws.Range("C5").Value = 4.9165802001953125
$ws.Range("D5").Value = "Hug"
$ws.Range("E5").Value = "None"

# Row 6 - new record
$ws.Range("A6").Value = 1.046940803527832
$ws.Range("B6").Value = "Distress"
$ws.Range("C6").Value = 4.1016740798950195
$ws.Range("D6").Value = "Discussion"
$ws.Range("E6").Value = "None"

# Row 7 - new record
$ws.Range("A7").Value = 1.7211663722991943
$ws.Range("B7").Value = "Joy"
$ws.Range("C7").Value = 2.2095818519592285
$ws.Range("D7").Value = "Congrat"
$ws.Range("E7").Value = "None"

# Row 8 - new record
$ws.Range("A8").Value = 0.7149765491485596
$ws.Range("B8").Value = "Distress"
$ws.Range("C8").Value = 3.1760647296905518
$ws.Range("D8").Value = "Bye"
$ws.Range("E8").Value = "Attention Deployment"

# Row 9 - new record
$ws.Range("A9").Value = 0
$ws.Range("B9").Value = "Hate"
$ws.Range("C9").Value = 1.3174933195114136
$ws.Range("D9").Value = "Fired"
$ws.Range("E9").Value = "Situation Modification"

# Row 10 - new record
$ws.Range("A10").Value = 0
$ws.Range("B10").Value = "Hate"
$ws.Range("C10").Value = 1.5745996236801147
$ws.Range("D10").Value = "Crash"
$ws.Range("E10").Value = "Cognitive Change"

# Row 11 - new record (this used to hold the first "PERSONALITY TRAITS" entry in col F)
$ws.Range("A11").Value = 2.372433662414551
$ws.Range("B11").Value = "Joy"
$ws.Range("C11").Value = 7.63873291015625
$ws.Range("D11").Value = "Profits"
$ws.Range("E11").Value = "None"

# The "PERSONALITY TRAITS" list (column F) moves from rows 3-7 down to rows 12-16
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("F7").ClearContents()

$ws.Range("F12").Value = "Low Conscientiousness"
$ws.Range("F13").Value = "Low Extraversion"
$ws.Range("F14").Value = "Low Neuroticism"
$ws.Range("F15").Value = "Low Agreeableness"
$ws.Range("F16").Value = "High Openness"

# The "STRATEGIES RELATED" list (column G) moves from rows 8-10 down to rows 17-21
# and gains two new entries at the top
$ws.Range("G8").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("G10").ClearContents()

$ws.Range("G17").Value = "[Situation Selection, Weakly]"
$ws.Range("G18").Value = "[Situation Modification, Strongly]"
$ws.Range("G19").Value = "[Attention Deployment, Strongly]"
$ws.Range("G20").Value = "[Cognitive Change, Strongly]"
$ws.Range("G21").Value = "[Response Modulation, Weakly]"

# The "DOMINANT PERSONALITY" entry (column H) moves from row 11 down to row 22
$ws.Range("H11").ClearContents()
$ws.Range("H22").Value = "Openness"
